# Actualización automática 2025-11-17 17:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M25").Value = 23307.46
$wsGrupo.Range("M26").Value = "4 de 24"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F25").Value = 23307.46
$wsMensual.Range("F26").Value = 28802.69

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 28359.25
$wsCumplimiento.Range("E12").Value = 6341.75
$wsCumplimiento.Range("F12").Value = 0.8172459006945045
$wsCumplimiento.Range("D14").Value = 28802.69
$wsCumplimiento.Range("E14").Value = 11975.05058948192
$wsCumplimiento.Range("F14").Value = 0.7063336414335146
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.15
